# Fruta / hortaliza, semanal
# Insert one new week of "Membrillo" price data (2 quality-grade rows) at the
# top of the historical table, pushing the existing rows (old 51-92) down to
# (53-94). This mirrors the existing weekly-record layout used throughout the
# sheet (two rows per date: "Primera" and "Segunda" quality grades).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current first record of this block
# (row 51). Excel shifts everything below down by two rows and the new rows
# inherit the formatting of the row that was previously there (keeps the
# date-formatted style on column D).
$ws.Rows("51:52").Insert()

# New row 51: "Primera" quality grade
$ws.Cells.Item(51, 1).Value = 8
$ws.Cells.Item(51, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(51, 3).Value = "Coquimbo"
$ws.Cells.Item(51, 4).Value = 45072
$ws.Cells.Item(51, 5).Value = 4
$ws.Cells.Item(51, 6).Value = "Fruta"
$ws.Cells.Item(51, 7).Value = 100104
$ws.Cells.Item(51, 8).Value = "Frutos de pepita"
$ws.Cells.Item(51, 9).Value = 100104003
$ws.Cells.Item(51, 10).Value = "Membrillo"
$ws.Cells.Item(51, 11).Value = "Champion"
$ws.Cells.Item(51, 12).Value = "Primera"
$ws.Cells.Item(51, 13).Value = 16
$ws.Cells.Item(51, 14).Value = 240000
$ws.Cells.Item(51, 15).Value = 250000
$ws.Cells.Item(51, 16).Value = 245000
$ws.Cells.Item(51, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(51, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(51, 19).Value = 544
$ws.Cells.Item(51, 20).Value = 450

# New row 52: "Segunda" quality grade
$ws.Cells.Item(52, 1).Value = 8
$ws.Cells.Item(52, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(52, 3).Value = "Coquimbo"
$ws.Cells.Item(52, 4).Value = 45072
$ws.Cells.Item(52, 5).Value = 4
$ws.Cells.Item(52, 6).Value = "Fruta"
$ws.Cells.Item(52, 7).Value = 100104
$ws.Cells.Item(52, 8).Value = "Frutos de pepita"
$ws.Cells.Item(52, 9).Value = 100104003
$ws.Cells.Item(52, 10).Value = "Membrillo"
$ws.Cells.Item(52, 11).Value = "Champion"
$ws.Cells.Item(52, 12).Value = "Segunda"
$ws.Cells.Item(52, 13).Value = 16
$ws.Cells.Item(52, 14).Value = 200000
$ws.Cells.Item(52, 15).Value = 210000
$ws.Cells.Item(52, 16).Value = 205000
$ws.Cells.Item(52, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(52, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(52, 19).Value = 456
$ws.Cells.Item(52, 20).Value = 450
